{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// A few movie titles were originally typed as two runs split around a\n// Word spell-check \"proofErr\" marker (e.g. \"Super \" + \"man\"). Rewrite\n// each paragraph's content as plain text so it collapses back into a\n// single clean run with no leftover proofErr markup.\nconst fixes = {\n  \"Super man\": \"Super man\",\n  \"Star wars\": \"Star wars\",\n  \"Harry potter\": \"Harry potter\"\n};\n\nfor (const paragraph of paragraphs.items) {\n  if (Object.prototype.hasOwnProperty.call(fixes, paragraph.text)) {\n    const range = paragraph.getRange();\n    range.clear();\n    await context.sync();\n    range.insertText(fixes[paragraph.text], \"Start\");\n  }\n}\n\n// The trailing empty paragraph becomes the new movie entry.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst lastRange = lastParagraph.getRange();\nlastRange.clear();\nawait context.sync();\nlastRange.insertText(\"Barbie vida de sereia\", \"Start\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Word's live spell-checker had split a few titles into two runs around a\n# <w:proofErr> pair (e.g. \"Super \" + \"man\"). Just overwriting Range.Text\n# only rewrites the first run and leaves the trailing run/proofErr behind,\n# so instead we delete the whole paragraph (content + its paragraph mark)\n# and insert a brand-new paragraph in its place holding a single clean\n# run with the full title. This removes the stray <w:proofErr/> markers\n# entirely, matching how Word itself collapses the runs once the\n# misspelling no longer exists.\nfunction Set-CleanParagraphText($index, $text) {\n    $paragraph = $d.Paragraphs.Item($index)\n    $paragraph.Range.Delete()\n    $following = $d.Paragraphs.Item($index)\n    $following.Range.InsertParagraphBefore()\n    $d.Paragraphs.Item($index).Range.Text = $text\n}\n\nSet-CleanParagraphText 2 \"Super man\"\nSet-CleanParagraphText 4 \"Star wars\"\nSet-CleanParagraphText 5 \"Harry potter\"\n\n# The trailing blank paragraph becomes the newly added movie.\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastParagraph.Range.Text = \"Barbie vida de sereia\"\n"}
